$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(-1.0144944190979, -1.167210817337036, 0.6551529765129089),
    @(-0.3874412775039673, -0.7050912380218506, 0.0914770737290382),
    @(-0.1458440721035003, -0.3762930035591125, -0.0704022198915481),
    @(0.2157881408929824, 0.3240640163421631, 0.0951422601938247),
    @(0.0403171069920063, 0.1484402567148208, -0.0852157026529312),
    @(0.1014036312699318, 0.3179553747177124, 0.0390953756868839),
    @(0.052381694316864, 0.1099557429552078, 0.0681114718317985),
    @(0.0522289797663688, -0.4196644127368927, 0.2273945808410644),
    @(0.0937678143382072, -0.1565342247486114, 0.0675006061792373),
    @(-0.0591012127697467, 0.0331394411623477, 0.0291688162833452)
)

$startRow = 22
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
    $ws.Cells.Item($row, 3).Value = $data[$i][2]
}
